$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 675
$ws.Range("I18").Value = 675
$ws.Range("K18").Value = 675
$ws.Range("M18").Value = -391
$ws.Range("H129").Value = 816.9677
$ws.Range("I129").Value = 412.35
$ws.Range("J129").Value = 1552.6364
$ws.Range("K129").Value = 1237.05
$ws.Range("L129").Value = 4657.9092
$ws.Range("M129").Value = 3762.95
$ws.Range("N129").Value = -14657.9092
$ws.Range("H137").Value = 2587263.5
$ws.Range("I137").Value = 1163828
$ws.Range("J137").Value = 6667778
$ws.Range("K137").Value = 3491484
$ws.Range("L137").Value = 20003334
$ws.Range("M137").Value = -3488934
$ws.Range("N137").Value = -20008434
$ws.Range("H141").Value = 1353.3729
$ws.Range("I141").Value = 943.9773
$ws.Range("J141").Value = 2554.2666
$ws.Range("K141").Value = 2831.9319
$ws.Range("L141").Value = 7662.7998
$ws.Range("M141").Value = 2348.0681
$ws.Range("N141").Value = -18022.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2665940
$ws.Range("I32").Value = 3708.7952
$ws.Range("J32").Value = 22753686
$ws.Range("K32").Value = 3708.7952
$ws.Range("L32").Value = 22753686
$ws.Range("M32").Value = -3421.7952
$ws.Range("N32").Value = -22754260
$ws.Range("H61").Value = 1276.9678
$ws.Range("I61").Value = 1403.1305
$ws.Range("J61").Value = 914.25
$ws.Range("K61").Value = 1403.1305
$ws.Range("L61").Value = 914.25
$ws.Range("M61").Value = -1191.1305
$ws.Range("N61").Value = -1338.25
$ws.Range("H74").Value = 4762860.5
$ws.Range("I74").Value = 940.8158
$ws.Range("K74").Value = 940.8158
$ws.Range("M74").Value = -66.81579999999997
$ws.Range("H77").Value = 4762860.5
$ws.Range("I77").Value = 940.8158
$ws.Range("K77").Value = 4704.079
$ws.Range("M77").Value = -336.0789999999997
$ws.Range("H132").Value = 58931.438
$ws.Range("I132").Value = 70443.664
$ws.Range("J132").Value = 3672.7334
$ws.Range("K132").Value = 211330.992
$ws.Range("L132").Value = 11018.2002
$ws.Range("M132").Value = -208800.992
$ws.Range("N132").Value = -16078.2002
$ws.Range("H136").Value = 1276.9678
$ws.Range("I136").Value = 1403.1305
$ws.Range("J136").Value = 914.25
$ws.Range("K136").Value = 4209.3915
$ws.Range("L136").Value = 2742.75
$ws.Range("M136").Value = -1659.3915
$ws.Range("N136").Value = -7842.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20558
$ws.Range("J41").Value = 21197.5
$ws.Range("L41").Value = 21197.5
$ws.Range("N41").Value = -22053.5
$ws.Range("H52").Value = 26319.572
$ws.Range("J52").Value = 26739.5
$ws.Range("L52").Value = 26739.5
$ws.Range("N52").Value = -27327.5
$ws.Range("H129").Value = 49749.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49749.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49749.5
$ws.Range("N129").Value = -59749.5
$ws.Range("H130").Value = 49980
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 49980
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 49980
$ws.Range("N130").Value = -60020
$ws.Range("H131").Value = 33851.668
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 33851.668
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 33851.668
$ws.Range("N131").Value = -43931.668
$ws.Range("H132").Value = 1739.2632
$ws.Range("I132").Value = 1678.102
$ws.Range("J132").Value = 2113.875
$ws.Range("K132").Value = 5034.306
$ws.Range("L132").Value = 6341.625
$ws.Range("M132").Value = -2504.306
$ws.Range("N132").Value = -11401.625
$ws.Range("H133").Value = 30141.25
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 30141.25
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 30141.25
$ws.Range("N133").Value = -35201.25
$ws.Range("H134").Value = 7313.5654
$ws.Range("I134").Value = 8841.117
$ws.Range("J134").Value = 2985.5
$ws.Range("K134").Value = 26523.351
$ws.Range("L134").Value = 8956.5
$ws.Range("M134").Value = -23988.351
$ws.Range("N134").Value = -14026.5
$ws.Range("H135").Value = 74999.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 74999.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 74999.75
$ws.Range("N135").Value = -85139.75
$ws.Range("H137").Value = 45000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 45000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -55200
$ws.Range("H138").Value = 46107
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 46107
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 46107
$ws.Range("N138").Value = -56387
$ws.Range("H139").Value = 53999.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 53999.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 53999.5
$ws.Range("N139").Value = -64279.5
$ws.Range("H140").Value = 58133.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 58133.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 58133.332
$ws.Range("N140").Value = -68493.33199999999
$ws.Range("H141").Value = 33462.8
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 33462.8
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 33462.8
$ws.Range("N141").Value = -43822.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1642.6578
$ws.Range("J39").Value = 1670.8379
$ws.Range("L39").Value = 5012.5137
$ws.Range("N39").Value = -5600.5137
$ws.Range("H80").Value = 8742.857
$ws.Range("J80").Value = 8742.857
$ws.Range("L80").Value = 26228.571
$ws.Range("N80").Value = -28100.571
$ws.Range("H83").Value = 8742.857
$ws.Range("J83").Value = 8742.857
$ws.Range("L83").Value = 78685.713
$ws.Range("N83").Value = -88045.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1623.3135
$ws.Range("I132").Value = 1292.7142
$ws.Range("J132").Value = 2178.72
$ws.Range("K132").Value = 3878.1426
$ws.Range("L132").Value = 6536.16
$ws.Range("M132").Value = -1348.1426
$ws.Range("N132").Value = -11596.16

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3383.6155
$ws.Range("I100").Value = 3548.3
$ws.Range("J100").Value = 2834.6667
$ws.Range("K100").Value = 3548.3
$ws.Range("L100").Value = 2834.6667
$ws.Range("M100").Value = -3007.3
$ws.Range("N100").Value = -3916.6667
$ws.Range("H132").Value = 1784
$ws.Range("I132").Value = 1496.5116
$ws.Range("J132").Value = 4874.5
$ws.Range("K132").Value = 4489.5348
$ws.Range("L132").Value = 14623.5
$ws.Range("M132").Value = -1959.5348
$ws.Range("N132").Value = -19683.5
$ws.Range("H133").Value = 41252.125
$ws.Range("J133").Value = 41252.125
$ws.Range("L133").Value = 41252.125
$ws.Range("N133").Value = -46312.125
$ws.Range("H136").Value = 1455.5363
$ws.Range("I136").Value = 1355.2
$ws.Range("J136").Value = 1849.7142
$ws.Range("K136").Value = 4065.6
$ws.Range("L136").Value = 5549.142599999999
$ws.Range("M136").Value = -1515.6
$ws.Range("N136").Value = -10649.1426

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 29900
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 29900
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 29900
$ws.Range("N119").Value = -39576
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H121").Value = 27546
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 27546
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 27546
$ws.Range("N121").Value = -31040
$ws.Range("H122").Value = 17780.072
$ws.Range("I122").Value = 27239
$ws.Range("J122").Value = 5168.1665
$ws.Range("K122").Value = 81717
$ws.Range("L122").Value = 15504.4995
$ws.Range("M122").Value = -79267
$ws.Range("N122").Value = -20404.4995
$ws.Range("H123").Value = 15111.25
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 15111.25
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 15111.25
$ws.Range("N123").Value = -24911.25
$ws.Range("H124").Value = 51330
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 51330
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 51330
$ws.Range("N124").Value = -61150
$ws.Range("H125").Value = 35799.582
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 35799.582
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 35799.582
$ws.Range("N125").Value = -45639.582
$ws.Range("H126").Value = 2452.3447
$ws.Range("I126").Value = 1996.7368
$ws.Range("J126").Value = 3318
$ws.Range("K126").Value = 5990.2104
$ws.Range("L126").Value = 9954
$ws.Range("M126").Value = -3520.2104
$ws.Range("N126").Value = -14894
$ws.Range("H127").Value = 40689.285
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 40689.285
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 40689.285
$ws.Range("N127").Value = -50609.285
$ws.Range("H128").Value = 49132.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49132.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49132.5
$ws.Range("N128").Value = -59092.5
$ws.Range("H129").Value = 36392.25
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 36392.25
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 36392.25
$ws.Range("N129").Value = -46392.25
$ws.Range("H130").Value = 45000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 45000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 1086.5802
$ws.Range("I132").Value = 965.73914
$ws.Range("J132").Value = 1781.4166
$ws.Range("K132").Value = 2897.21742
$ws.Range("L132").Value = 5344.2498
$ws.Range("M132").Value = -367.2174199999999
$ws.Range("N132").Value = -10404.2498
$ws.Range("H133").Value = 35486
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 35486
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 35486
$ws.Range("N133").Value = -45606
$ws.Range("H135").Value = 49945.625
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 49945.625
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 49945.625
$ws.Range("N135").Value = -60085.625
$ws.Range("H136").Value = 1742.78
$ws.Range("I136").Value = 1768.2727
$ws.Range("J136").Value = 1555.8334
$ws.Range("K136").Value = 5304.8181
$ws.Range("L136").Value = 4667.5002
$ws.Range("M136").Value = -2754.8181
$ws.Range("N136").Value = -9767.5002
$ws.Range("H137").Value = 27995
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 27995
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 27995
$ws.Range("N137").Value = -38195
$ws.Range("H138").Value = 40426.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 40426.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 40426.5
$ws.Range("N138").Value = -50706.5
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 28359.428
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 28359.428
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 28359.428
$ws.Range("N140").Value = -38719.428
$ws.Range("H141").Value = 63691.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 63691.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 63691.25
$ws.Range("N141").Value = -74051.25
